# Finished Week 13 logging
# Update row 3 ("R") Target Depth data on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# OFF sheet
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 169
$wsOff.Range("C3").Value = 111
$wsOff.Range("D3").Value = 42
$wsOff.Range("E3").Value = 21
$wsOff.Range("F3").Value = 3

# DEF sheet
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 159
$wsDef.Range("C3").Value = 97
$wsDef.Range("D3").Value = 49
$wsDef.Range("E3").Value = 29
$wsDef.Range("F3").Value = 2
